$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "Trening" header ---
# Copy header style from E1 so F1 matches the other header cells (bold, border, centered)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# --- Column A: convert text dates to real date/time values, rows 2-5 ---
# (the commit registers both an intermediate lower-case format and the
# final upper-case format that actually gets applied)
$ws.Cells.Item(2, 1).Value = 45684
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(3, 1).Value = 45684
$ws.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(4, 1).Value = 45684
$ws.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(5, 1).Value = 45684
$ws.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 2 & 3: label the existing training-block rows ("Duża Gra") ---
$ws.Range("F2").Value = "Duża Gra"
$ws.Range("F3").Value = "Duża Gra"

# --- Row 4 & 5: new rows for the second training block ("Mała Gra") ---
$ws.Range("E4").Value = "10-15"
$ws.Range("F4").Value = "Mała Gra"

$ws.Range("E5").Value = "5-10"
$ws.Range("F5").Value = "Mała Gra"

Write-Output "done"
